# "Folhas para os dois mapas" - split the single sheet into two sheets,
# one per map (Mapa2a / Mapa2b), duplicating the original data.

$wb = $excel.ActiveWorkbook

# Original (only) sheet currently holds the "Mapa2a" data.
$original = $wb.Worksheets.Item("Folha1")

# Duplicate it, inserting the copy right before the original so the copy
# becomes the first tab (gets a fresh, higher sheetId; the original keeps
# its existing sheetId but moves to second position).
$original.Copy($original)

# After the copy, tab order is: [copy, original]
$mapa2a = $wb.Worksheets.Item(1)
$mapa2b = $wb.Worksheets.Item(2)

$mapa2a.Name = "Mapa2a"
$mapa2b.Name = "Mapa2b"

# The duplicate still says "Mapa2a" in A1 - relabel the second sheet.
$mapa2b.Range("A1").Value = "Mapa2b"

# Restore the original's lingering selection state and reselect the first
# (now Mapa2a) tab as the active sheet/view.
$mapa2b.Range("D18").Select()
$mapa2a.Select()
